$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 368.42856
$ws.Range("I33").Value = 368.42856
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 368.42856
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -139.42856

$ws.Range("H41").Value = 449.7143
$ws.Range("J41").Value = 195.33333
$ws.Range("L41").Value = 195.33333
$ws.Range("N41").Value = -1075.33333

$ws.Range("H86").Value = 2120.7827
$ws.Range("I86").Value = 2571.3845
$ws.Range("J86").Value = 1535
$ws.Range("K86").Value = 2571.3845
$ws.Range("L86").Value = 1535
$ws.Range("M86").Value = -1448.3845
$ws.Range("N86").Value = -3781

$ws.Range("H89").Value = 2120.7827
$ws.Range("I89").Value = 2571.3845
$ws.Range("J89").Value = 1535
$ws.Range("K89").Value = 12856.9225
$ws.Range("L89").Value = 7675
$ws.Range("M89").Value = -7240.922500000001
$ws.Range("N89").Value = -18907

$ws.Range("H98").Value = 2802.5
$ws.Range("I98").Value = 2802.5
$ws.Range("K98").Value = 2802.5
$ws.Range("M98").Value = -1304.5

$ws.Range("H107").Value = 549.9286
$ws.Range("I107").Value = 530.6923
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 530.6923
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 1389.3077
$ws.Range("N107").Value = -4640

$ws.Range("H116").Value = 5709.4287
$ws.Range("I116").Value = 5364.074
$ws.Range("K116").Value = 5364.074
$ws.Range("M116").Value = -1922.074

$ws.Range("H122").Value = 2802.5
$ws.Range("I122").Value = 2802.5
$ws.Range("K122").Value = 8407.5
$ws.Range("M122").Value = -5957.5

$ws.Range("H125").Value = 3182.3462
$ws.Range("I125").Value = 2897.2307
$ws.Range("J125").Value = 3467.4614
$ws.Range("K125").Value = 26075.0763
$ws.Range("L125").Value = 31207.1526
$ws.Range("M125").Value = -23615.0763
$ws.Range("N125").Value = -36127.1526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5217.5
$ws.Range("I32").Value = 5921.2
$ws.Range("J32").Value = 2480.889
$ws.Range("K32").Value = 5921.2
$ws.Range("L32").Value = 2480.889
$ws.Range("M32").Value = -5634.2
$ws.Range("N32").Value = -3054.889

$ws.Range("H45").Value = 1224.6666
$ws.Range("I45").Value = 1169.6
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1169.6
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -792.5999999999999
$ws.Range("N45").Value = -2254

$ws.Range("H98").Value = 15998
$ws.Range("J98").Value = 15998
$ws.Range("L98").Value = 15998
$ws.Range("N98").Value = -21988

$ws.Range("H103").Value = 48000
$ws.Range("J103").Value = 48000
$ws.Range("L103").Value = 48000
$ws.Range("N103").Value = -50344

$ws.Range("H104").Value = 47750
$ws.Range("J104").Value = 47750
$ws.Range("L104").Value = 47750
$ws.Range("N104").Value = -54738

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0

$ws.Range("H110").Value = 8356.5
$ws.Range("I110").Value = 4700
$ws.Range("J110").Value = 12013
$ws.Range("K110").Value = 4700
$ws.Range("L110").Value = 12013
$ws.Range("M110").Value = -2655
$ws.Range("N110").Value = -16103

$ws.Range("H115").Value = 48000
$ws.Range("J115").Value = 48000
$ws.Range("L115").Value = 48000
$ws.Range("N115").Value = -51134

$ws.Range("H119").Value = 10698
$ws.Range("J119").Value = 10698
$ws.Range("L119").Value = 10698
$ws.Range("N119").Value = -20374

$ws.Range("H122").Value = 2405.6
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 2257
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 6771
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -11671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 42880.5
$ws.Range("J120").Value = 42880.5
$ws.Range("L120").Value = 42880.5
$ws.Range("N120").Value = -52556.5

$ws.Range("H130").Value = 47780
$ws.Range("J130").Value = 47780
$ws.Range("L130").Value = 47780
$ws.Range("N130").Value = -57820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 337010
$ws.Range("I3").Value = 337010
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1011030
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -1010918

$ws.Range("H68").Value = 2151.4443
$ws.Range("I68").Value = 1251.5
$ws.Range("J68").Value = 3951.3333
$ws.Range("K68").Value = 3754.5
$ws.Range("L68").Value = 11853.9999
$ws.Range("M68").Value = -2943.5
$ws.Range("N68").Value = -13475.9999

$ws.Range("H71").Value = 2151.4443
$ws.Range("I71").Value = 1251.5
$ws.Range("J71").Value = 3951.3333
$ws.Range("K71").Value = 11263.5
$ws.Range("L71").Value = 35561.9997
$ws.Range("M71").Value = -7207.5
$ws.Range("N71").Value = -43673.9997

$ws.Range("H87").Value = 10000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 10000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H122").Value = 1444.6154
$ws.Range("J122").Value = 1820
$ws.Range("L122").Value = 16380
$ws.Range("N122").Value = -21280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 26330
$ws.Range("J95").Value = 26330
$ws.Range("L95").Value = 26330
$ws.Range("N95").Value = -31822

$ws.Range("H107").Value = 235
$ws.Range("I107").Value = 246.66667
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 246.66667
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 1673.33333
$ws.Range("N107").Value = -4040

$ws.Range("H113").Value = 4690.579
$ws.Range("I113").Value = 4251.5
$ws.Range("J113").Value = 5920
$ws.Range("K113").Value = 4251.5
$ws.Range("L113").Value = 5920
$ws.Range("M113").Value = -2081.5
$ws.Range("N113").Value = -10260

$ws.Range("H119").Value = 48000
$ws.Range("J119").Value = 48000
$ws.Range("L119").Value = 48000
$ws.Range("N119").Value = -57676

$ws.Range("H122").Value = 1587.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1587.5
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").Value = 4762.5
$ws.Range("N122").Value = -9662.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3070.5
$ws.Range("I61").Value = 2434
$ws.Range("J61").Value = 4980
$ws.Range("K61").Value = 2434
$ws.Range("L61").Value = 4980
$ws.Range("M61").Value = -2232
$ws.Range("N61").Value = -5384

$ws.Range("H113").Value = 3070.5
$ws.Range("I113").Value = 2434
$ws.Range("J113").Value = 4980
$ws.Range("K113").Value = 2434
$ws.Range("L113").Value = 4980
$ws.Range("M113").Value = -264
$ws.Range("N113").Value = -9320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 33000
$ws.Range("J98").Value = 33000
$ws.Range("L98").Value = 33000
$ws.Range("N98").Value = -38990

$ws.Range("H107").Value = 470.4
$ws.Range("I107").Value = 470.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1411.2
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 508.8000000000002

$ws.Range("H113").Value = 91801.27
$ws.Range("I113").Value = 111701.555
$ws.Range("K113").Value = 335104.665
$ws.Range("M113").Value = -332934.665

$ws.Range("H119").Value = 25560
$ws.Range("J119").Value = 25560
$ws.Range("L119").Value = 25560
$ws.Range("N119").Value = -35236

$ws.Range("H126").Value = 935.4286
$ws.Range("I126").Value = 703.38464
$ws.Range("J126").Value = 1312.5
$ws.Range("K126").Value = 2110.15392
$ws.Range("L126").Value = 3937.5
$ws.Range("M126").Value = 359.8460800000003
$ws.Range("N126").Value = -8877.5
